$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 61

$ws.Range("A$row").Value = 60
$ws.Range("B$row").Value = "india"
$ws.Range("C$row").Value = "isl"
$ws.Range("D$row").Value = "2023-2024"
$ws.Range("E$row").Value = 45282.64583333334
$ws.Range("F$row").Value = "East Bengal"
$ws.Range("G$row").Value = 0
$ws.Range("H$row").Value = "Odisha FC"
$ws.Range("I$row").Value = 0
$ws.Range("J$row").Value = 2.37
$ws.Range("K$row").Value = "17/12/2023 15:42"
$ws.Range("L$row").Value = 2.84
$ws.Range("M$row").Value = "22/12/2023 15:29"
$ws.Range("N$row").Value = 3.31
$ws.Range("O$row").Value = "17/12/2023 15:42"
$ws.Range("P$row").Value = 3.4
$ws.Range("Q$row").Value = "22/12/2023 15:29"
$ws.Range("R$row").Value = 2.86
$ws.Range("S$row").Value = "17/12/2023 15:42"
$ws.Range("T$row").Value = 2.52
$ws.Range("U$row").Value = "22/12/2023 15:29"
$ws.Range("V$row").Value = "https://www.betexplorer.com/football/india/isl/east-bengal-odisha-fc/ObkJrnft/"

# Copy cell formatting (style) from the row above to match the sheet's
# established pattern (bold/bordered index column, date-formatted match date).
$ws.Range("A60").Copy()
$ws.Range("A$row").PasteSpecial(-4122)
$ws.Range("E60").Copy()
$ws.Range("E$row").PasteSpecial(-4122)

$excel.CutCopyMode = 0
